$wb = $excel.ActiveWorkbook

# Select all cells on the Germany sheet (matches the recorded selection
# state in the target: sqref="A1:XFD1048576", no tabSelected).
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Cells.Select()

# New "Swiss" sheet: same layout/styles/col-widths as "Belgium", so copy
# it and then tweak the market-specific cells.
$belgium = $wb.Worksheets.Item("Belgium")
$czech = $wb.Worksheets.Item("Czech")
$belgium.Copy($null, $czech)

$swiss = $wb.Worksheets.Item(4)
$swiss.Name = "Swiss"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2344"

$swiss.Activate()
$swiss.Range("B10").Select()
